# Update of league bases, 30-05-2024 23:16
# 1) Header rename: ht_goals_h -> HTHG, ht_goals_a -> HTAG
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "HTHG"
$ws.Range("I1").Value = "HTAG"

# 2) Row data re-shuffle (columns B..AD only; column A/id is left untouched).
#    Snapshot every affected row *before* writing anything back, since some
#    groups are cyclic permutations (not simple pairwise swaps).
$rows = @(93, 94, 95, 96, 97, 98, 102, 103, 114, 115, 135, 136, 157, 158)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = $ws.Range("B$r`:AD$r").Value2
}

# destination row -> source row (source data, as it existed before edits,
# that should end up in the destination row)
$mapping = @{
    93  = 95
    94  = 98
    95  = 94
    96  = 97
    97  = 96
    98  = 93
    102 = 103
    103 = 102
    114 = 115
    115 = 114
    135 = 136
    136 = 135
    157 = 158
    158 = 157
}

foreach ($dest in $rows) {
    $src = $mapping[$dest]
    $ws.Range("B$dest`:AD$dest").Value = $snapshot[$src]
}
